$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "307"
$t.Cell(5,1).Range.Text  = "0.00003"
$t.Cell(6,1).Range.Text  = "0.00049"
$t.Cell(7,1).Range.Text  = "0.00017"
$t.Cell(8,1).Range.Text  = "0.00003"
$t.Cell(9,1).Range.Text  = "0.00038"
$t.Cell(10,1).Range.Text = "0.00040"
$t.Cell(11,1).Range.Text = "0.00042"
$t.Cell(12,1).Range.Text = "0.06525"

$t.Cell(44,1).Range.Text = "99.7"
$t.Cell(45,1).Range.Text = "0.07"
$t.Cell(46,1).Range.Text = "22"
